$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.114.21'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '2.423.34'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").Value = '2.422.97'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("E12").Value = '  -3.48%  '
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000174'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.09%  '
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").Value = '62.058.90'
$ws.Range("E17").Value = '  +0.15%  '
$ws.Range("D18").Value = '2.438.50'
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.59%  '
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("E22").Value = '  -1.36%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.65%  '
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '560.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.21%  '
$ws.Range("D28").Value = '2.544.35'
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -1.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("E32").Value = '  -5.54%  '
$ws.Range("E33").Value = '  -1.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.86'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.42%  '
$ws.Range("E35").Value = '  -3.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  -1.42%  '
$ws.Range("E38").Value = '  -1.17%  '
$ws.Range("E39").Value = '  -4.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '152.23'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("E44").Value = '  -2.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.89'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0531'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.596'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("E51").Value = '  -0.34%  '
